$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.185.74"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "'1.853.18"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'235.36"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4773"
$ws.Range("E7").Value = "  -2.95%  "
$ws.Range("D8").Value = "'0.2808"
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "'0.06458"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "'1.851.22"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").Value = "'0.07364"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "'16.20"
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").Value = "'5.085"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'87.01"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.6447"
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("D16").Value = "'30.144.57"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'13.11"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("D19").Value = "'0.000007527"
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").Value = "'226.20"
$ws.Range("E20").Value = "  +17.26%  "
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.274"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "'2.086.82"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'6.069"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'9.197"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D26").Value = "'163.12"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'18.49"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "'1.916"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").Value = "'1.439"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").Value = "'0.09178"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'4.227"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "'3.940"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").Value = "'0.04965"
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.138"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'0.7253"
$ws.Range("E35").Value = "  -1.86%  "
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").Value = "'0.01834"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("D39").Value = "'0.8971"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").Value = "'2.028"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "'5.928"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "'105.77"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.4229"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").Value = "'7.343"
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("D46").Value = "'0.1309"
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("D47").Value = "'63.94"
$ws.Range("E47").Value = "  -7.02%  "
$ws.Range("D48").Value = "'1.501"
$ws.Range("E48").Value = "  +7.29%  "
$ws.Range("D49").Value = "'8.680"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").Value = "'33.72"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").Value = "'0.05652"
$ws.Range("E51").Value = "  -3.43%  "
